$wb = $excel.ActiveWorkbook

# The localization status changed from "Ready for handoff" to "In Translation".
# That shared string is used in four cells across the three sheets:
#   Overview!E2 (zh-cn status), Overview!F2 (de-de status),
#   zh-cn!C2 (Status column) and de-de!C2 (Status column).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The "Status" columns were narrowed to fit the new, shorter text
# (Overview columns E & F, and column C on the zh-cn / de-de sheets).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
